$wb = $excel.ActiveWorkbook

# --- Sheet: "Range Status" ---
# Zero out the Species (no.) column (B2:B7) and remove the Species (perc.) column (C2:C7)
$wsRange = $wb.Worksheets.Item("Range Status")
$wsRange.Range("B2:B7").Value = 0
$wsRange.Range("C2:C7").ClearContents()

# --- Sheet: "Species qualification" ---
# "Range Analysis" row (row 5) selected-for-analysis count drops to 0
$wsSpecies = $wb.Worksheets.Item("Species qualification")
$wsSpecies.Range("B5").Value = 0

# --- Sheet: "High Priority break-up" ---
# The "Range" row disappears entirely and the "IUCN" row absorbs all species,
# becoming the sole remaining row.
$wsBreakup = $wb.Worksheets.Item("High Priority break-up")
$wsBreakup.Range("A2").Value = "IUCN"
$wsBreakup.Range("B2").Value = 13
$wsBreakup.Range("C2").Value = 100
$wsBreakup.Range("D2").Value = 13
$wsBreakup.Range("E2").Value = 100
$wsBreakup.Rows.Item(3).Delete()
